$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Septiembre de 2020 a las 02:38"

# --- Row 4: Estados Unidos - refreshed totals ---
$ws.Range("B4").Value = 7184849
$ws.Range("C4").Value = 44733
$ws.Range("D4").Value = 4431052
$ws.Range("E4").Value = 2546281
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 920
$ws.Range("H4").Value = 207516

# --- Row 29 - refreshed totals ---
$ws.Range("B29").Value = 149094
$ws.Range("C29").Value = 1341
$ws.Range("D29").Value = 128707
$ws.Range("E29").Value = 11138
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 9249

# --- Rows 53/54: Venezuela overtakes Costa Rica in ranking ---
# Row 53 now becomes Venezuela with fresh data
$ws.Range("A53").Value = "Venezuela"
$ws.Range("B53").Value = 70406
$ws.Range("C53").Value = 967
$ws.Range("D53").Value = 59745
$ws.Range("E53").Value = 10080
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 7
$ws.Range("H53").Value = 581

# Row 54 now becomes Costa Rica, keeping the previous Costa Rica data
$ws.Range("A54").Value = "Costa Rica"
$ws.Range("B54").Value = 69459
$ws.Range("C54").Value = 1400
$ws.Range("D54").Value = 26554
$ws.Range("E54").Value = 42110
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 14
$ws.Range("H54").Value = 795

# --- Row 72 - refreshed totals ---
$ws.Range("B72").Value = 36404
$ws.Range("C72").Value = 833
$ws.Range("D72").Value = 20502
$ws.Range("E72").Value = 15159
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 16
$ws.Range("H72").Value = 743

# --- Row 128 - refreshed totals ---
$ws.Range("B128").Value = 4789
$ws.Range("C128").Value = 10
$ws.Range("D128").Value = 4578
$ws.Range("E128").Value = 109
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 102

# --- Row 160 - refreshed totals ---
$ws.Range("B160").Value = 1663
$ws.Range("C160").Value = 9
$ws.Range("D160").Value = 1369
$ws.Range("E160").Value = 272
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 22

# --- Row 170 - refreshed totals ---
$ws.Range("B170").Value = 727
$ws.Range("C170").Value = 4
$ws.Range("D170").Value = 676
$ws.Range("E170").Value = 9
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 42

# --- Row 177 - refreshed totals ---
$ws.Range("B177").Value = 477
$ws.Range("C177").Value = 1
$ws.Range("D177").Value = 462
$ws.Range("E177").Value = 14
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 1

# --- Rows 215/216: Montserrat overtakes Islas Malvinas in ranking ---
# Row 215 now becomes Montserrat (keeping the previous Montserrat data)
$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1

# Row 216 now becomes Islas Malvinas (keeping the previous Islas Malvinas data)
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0
